$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the touched Price/Volume cells stay text (values contain
# multi-dot thousand separators and fixed-width percent strings that
# Excel would otherwise coerce into numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.476.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.847.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "265.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5214"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3286"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06822"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7792"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07760"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.837.11"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.026"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9987"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.98"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007987"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9989"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.529.98"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.078.70"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.649"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.578"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.016"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.42"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.188"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -8.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.661"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.15"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.189"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.146"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08753"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7261"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.137"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.106"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01788"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.233"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4922"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9157"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "111.39"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.093"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.778"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4194"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05942"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.122"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1247"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.06"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8890"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.24%  "
